$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.08513380217814247
$ws.Range("E2").Value = 0.08513380217814247

# Row 3
$ws.Range("D3").Value = 0.4205986492282449
$ws.Range("E3").Value = 0.4205986492282449

# Row 4
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.8194362898292942
$ws.Range("E4").Value = 0.8194362898292942

# Row 5
$ws.Range("D5").Value = 0.4904466957091531
$ws.Range("E5").Value = 0.4904466957091531

# Row 6
$ws.Range("D6").Value = 0.1085483803470201
$ws.Range("E6").Value = 0.1085483803470201

# Row 7
$ws.Range("D7").Value = 0.07459593390281857
$ws.Range("E7").Value = 0.9254040660971814

# Row 8
$ws.Range("D8").Value = 0.4277093061739183
$ws.Range("E8").Value = 0.5722906938260817

# Row 9
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = 0.6682604480051199
$ws.Range("E9").Value = 0.3317395519948801

# Row 10
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.5705164301900247
$ws.Range("E10").Value = 0.4294835698099753

# Row 11
$ws.Range("D11").Value = 0.3155920684447934
$ws.Range("E11").Value = 0.6844079315552066
$ws.Range("F11").Value = 0.869810938835144
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.08391905908028191
$ws.Range("E12").Value = 0.08391905908028191

# Row 13
$ws.Range("D13").Value = 0.1810007622017242
$ws.Range("E13").Value = 0.1810007622017242

# Row 14
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = 0.9196271373364906
$ws.Range("E14").Value = 0.9196271373364906

# Row 15
$ws.Range("C15").Value = $false
$ws.Range("D15").Value = 0.6225750000391651
$ws.Range("E15").Value = 0.6225750000391651

# Row 16
$ws.Range("D16").Value = 0.03651048164369296
$ws.Range("E16").Value = 0.03651048164369296

# Row 17
$ws.Range("D17").Value = 0.1995863608717966
$ws.Range("E17").Value = 0.8004136391282034

# Row 18
$ws.Range("D18").Value = 0.4682641715584517
$ws.Range("E18").Value = 0.5317358284415483

# Row 19
$ws.Range("C19").Value = $true
$ws.Range("D19").Value = 0.796525482954875
$ws.Range("E19").Value = 0.203474517045125

# Row 20
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = 0.7103567116998867
$ws.Range("E20").Value = 0.2896432883001133

# Row 21
$ws.Range("D21").Value = 0.2333063837129989
$ws.Range("E21").Value = 0.7666936162870011
$ws.Range("F21").Value = 0.8215096592903137
$ws.Range("G21").Value = 0.6
